# article 90 is live
#
# Row 7 holds three "blog" card cells (C7, E7, I7) that reference a
# rotating "ser: NN" value. A new article (ser: 90) goes live, bumping
# the previous occupants down the chain:
#   C7 (ser: 89) -> ser: 90   (new article)
#   E7 (ser: 88) -> ser: 89   (was C7's value)
#   I7 (ser: 87) -> ser: 88   (was E7's value)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 90"
$ws.Range("E7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 89"
$ws.Range("I7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 88"
